$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above the current last data row (row 81) -------------
# Excel's row Insert copies formatting from the row pushed down, so the new
# row 81 picks up the date-time style (s="1") on column A automatically, and
# the old row 81 becomes row 82.
$ws.Rows.Item(81).Insert()

# New row 81: 2024-06-11, zero-volume bar.
$ws.Cells.Item(81, 1).Value2 = 45454.2916666667
$ws.Cells.Item(81, 2).Value2 = 0
$ws.Cells.Item(81, 3).Value2 = 6.26000022888184
$ws.Cells.Item(81, 4).Value2 = 6.26000022888184
$ws.Cells.Item(81, 5).Value2 = 6.26000022888184
$ws.Cells.Item(81, 6).Value2 = 6.26000022888184
# Column G stores the adj_close as TEXT (shared string), matching the
# source file's column F value formatted as a string. Force text entry via
# a temporary "@" number format, then drop back to the default "Normal"
# style so no stray s="..." attribute is left on the cell.
$ws.Cells.Item(81, 7).NumberFormat = "@"
$ws.Cells.Item(81, 7).Value = "6.26000022888184"
$ws.Cells.Item(81, 7).Style = "Normal"
$ws.Cells.Item(81, 8).Value = "PAL.MI"

# Row 82 keeps the values that used to live in row 81, except its date moves
# from 2024-06-12 14:28:34 (45455.6031712963) to 2024-06-12 (45455.2916666667).
$ws.Cells.Item(82, 1).Value2 = 45455.2916666667

# --- Append a brand-new row 83 ----------------------------------------------
# This row sits beyond the sheet's previous used range, so copy formatting
# from row 82 first (keeps column A on the shared date style s="1" without
# minting a new style entry), then overwrite the values.
$ws.Cells.Item(82, 1).Copy($ws.Cells.Item(83, 1))

$ws.Cells.Item(83, 1).Value2 = 45456.6193865741
$ws.Cells.Item(83, 2).Value2 = 1800
$ws.Cells.Item(83, 3).Value2 = 6.19999980926514
$ws.Cells.Item(83, 4).Value2 = 6.15999984741211
$ws.Cells.Item(83, 5).Value2 = 6.15999984741211
$ws.Cells.Item(83, 6).Value2 = 6.15999984741211
$ws.Cells.Item(83, 7).NumberFormat = "@"
$ws.Cells.Item(83, 7).Value = "6.15999984741211"
$ws.Cells.Item(83, 7).Style = "Normal"
$ws.Cells.Item(83, 8).Value = "PAL.MI"
